{"js": "// The source document has the hyperlink text \"https://www.youtube.com/watch?v=8ZPsZBcue50\"\n// split across three separate runs (\"https://www\", \".\", \"youtube.com/watch?v=8ZPsZBcue50\"),\n// and the timestamp text \"1:00:49\" is a typo for \"1:07:28\".\n//\n// This script:\n//   1) Re-types the hyperlink text so Word collapses it into a single contiguous run\n//      (the rendered/visible text and the hyperlink target are unchanged).\n//   2) Fixes the timestamp text \"1:00:49\" -> \"1:07:28\", preserving its run formatting\n//      (bold, blue color, Bookman Old Style font).\n\nconst body = context.document.body;\nconst searchOptions = { matchCase: true, matchWildcards: false };\n\n// --- 1) Normalize the split hyperlink text into one run -------------------------------\nconst hyperlinkText = \"https://www.youtube.com/watch?v=8ZPsZBcue50\";\n\nlet results = body.search(hyperlinkText, searchOptions);\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  let target = results.items[0];\n\n  // Re-insert the exact same text. A byte-for-byte identical replacement is treated as a\n  // no-op by the run model (it would leave the three original runs untouched), so first\n  // stamp a one-character sentinel to force Word to rebuild the range as a single run,\n  // then immediately replace that sentinel-bearing text with the final, correct text.\n  target.insertText(hyperlinkText + \"\\u0001\", Word.InsertLocation.replace);\n  await context.sync();\n\n  results = body.search(hyperlinkText + \"\\u0001\", searchOptions);\n  results.load(\"text\");\n  await context.sync();\n\n  target = results.items[0];\n  target.insertText(hyperlinkText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- 2) Fix the timestamp \"1:00:49\" -> \"1:07:28\" ---------------------------------------\nresults = body.search(\"1:00:49\", searchOptions);\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const target = results.items[0];\n  target.insertText(\"1:07:28\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# The source document has the hyperlink text \"https://www.youtube.com/watch?v=8ZPsZBcue50\"\n# split across three separate runs (\"https://www\", \".\", \"youtube.com/watch?v=8ZPsZBcue50\"),\n# and the timestamp text \"1:00:49\" is a typo for \"1:07:28\".\n#\n# This script:\n#   1) Re-types the hyperlink text so Word collapses it into a single contiguous run\n#      (the rendered/visible text and the hyperlink target are unchanged).\n#   2) Fixes the timestamp text \"1:00:49\" -> \"1:07:28\", preserving its run formatting\n#      (bold, blue color, Bookman Old Style font).\n\n$d = $word.ActiveDocument\n\n# --- 1) Normalize the split hyperlink text into one run ---------------------------------\n$hyperlinkText = \"https://www.youtube.com/watch?v=8ZPsZBcue50\"\n$sentinel = [char]1\n\n$range = $d.Content\n$find = $range.Find\n$find.Text = $hyperlinkText\n$find.MatchCase = $true\n$find.MatchWildcards = $false\n$found = $find.Execute()\n\nif ($found) {\n    # Re-insert the exact same text. A byte-for-byte identical replacement is treated as a\n    # no-op by the run model (it would leave the three original runs untouched), so first\n    # stamp a one-character sentinel to force Word to rebuild the range as a single run,\n    # then immediately replace that sentinel-bearing text with the final, correct text.\n    $range.Text = $hyperlinkText + $sentinel\n\n    $range2 = $d.Content\n    $find2 = $range2.Find\n    $find2.Text = $hyperlinkText + $sentinel\n    $find2.MatchCase = $true\n    $find2.MatchWildcards = $false\n    $found2 = $find2.Execute()\n    if ($found2) {\n        $range2.Text = $hyperlinkText\n    }\n}\n\n# --- 2) Fix the timestamp \"1:00:49\" -> \"1:07:28\" -----------------------------------------\n$tsRange = $d.Content\n$tsFind = $tsRange.Find\n$tsFind.Text = \"1:00:49\"\n$tsFind.MatchCase = $true\n$tsFind.MatchWildcards = $false\n$tsFound = $tsFind.Execute()\n\nif ($tsFound) {\n    $tsRange.Text = \"1:07:28\"\n}\n"}
